$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F, copying the header formatting from an
# existing header cell (bold font, border, centered alignment)
$ws.Range("B1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "time_taken"

# Timestamps for rows 2-21 (time_taken values)
$timestamps = @(
    "2021-10-05 13:41:21.573315",
    "2021-10-05 13:41:21.573327",
    "2021-10-05 13:41:21.573331",
    "2021-10-05 13:41:21.573334",
    "2021-10-05 13:41:21.573337",
    "2021-10-05 13:41:21.573340",
    "2021-10-05 13:41:21.573343",
    "2021-10-05 13:41:21.573347",
    "2021-10-05 13:41:21.573350",
    "2021-10-05 13:41:21.573353",
    "2021-10-05 13:41:21.573356",
    "2021-10-05 13:41:21.573359",
    "2021-10-05 13:41:21.573362",
    "2021-10-05 13:41:21.573365",
    "2021-10-05 13:41:21.573368",
    "2021-10-05 13:41:21.573371",
    "2021-10-05 13:41:21.573375",
    "2021-10-05 13:41:21.573378",
    "2021-10-05 13:41:21.573381",
    "2021-10-05 13:41:21.573384"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
